# Fruta / hortaliza, semanal
# Insert a new weekly price-record row for "Apio" (Vega Monumental Concepción)
# at row 75, pushing the existing rows 75..174 down to 76..175.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 75 (shifts rows 75:174 down to 76:175).
$ws.Rows.Item(75).Insert()

# Populate the newly inserted row with the new week's data.
$ws.Range("A75").Value = 11
$ws.Range("B75").Value = "Vega Monumental Concepción"
$ws.Range("C75").Value = "Bíobío"
$ws.Range("D75").Value = 44539
$ws.Range("D75").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E75").Value = 8
$ws.Range("F75").Value = 100112017
$ws.Range("G75").Value = "Apio"
$ws.Range("H75").Value = "Americana (o)"
$ws.Range("I75").Value = "Primera"
$ws.Range("J75").Value = 270
$ws.Range("K75").Value = 6000
$ws.Range("L75").Value = 6500
$ws.Range("M75").Value = 6278
$ws.Range("N75").Value = "$/docena de matas"
$ws.Range("O75").Value = "Región de Coquimbo"
$ws.Range("P75").Value = 1046
$ws.Range("Q75").Value = 6
$ws.Range("R75").Value = "Hortaliza"
